$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new task row ("Flash Memory") right before the old row 16
#    (INTEGRATION header). This shifts rows 16-19 down to 17-20, and Excel
#    automatically slides the mergeCell ref (B16:C16 -> B17:C17) for us.
# ---------------------------------------------------------------------------
$ws.Rows("16:16").Insert()

# The freshly inserted row copies formatting down from row 15 (Baton
# Sensing) -- clear that so only the cells we want to be populated keep a
# style.
$ws.Range("F16:K16").ClearContents()
$ws.Range("F16:K16").ClearFormats()

# New row content: "Flash Memory" task label, a green marker (I16) and a
# blue in-progress marker (J16), matching the other task rows.
$ws.Range("C16").Value = "Flash Memory"
$ws.Range("C5").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = "Flash Memory"

# ---------------------------------------------------------------------------
# 2. Re-colour a handful of existing Gantt bar cells.
#    The project re-used the stock "red" fill as a new "green" (theme
#    Accent 6) fill, and shuffled a few bars between blue / cleared / grey.
# ---------------------------------------------------------------------------

# I5, I7, I14, I15 (and the new I16) become the new green marker.
$ws.Range("I5").Interior.ThemeColor = 10
$ws.Range("I5").Interior.TintAndShade = 0

$ws.Range("I7").Interior.ThemeColor = 10
$ws.Range("I7").Interior.TintAndShade = 0

$ws.Range("I14").Interior.ThemeColor = 10
$ws.Range("I14").Interior.TintAndShade = 0

$ws.Range("I15").Interior.ThemeColor = 10
$ws.Range("I15").Interior.TintAndShade = 0

$ws.Range("I16").Interior.ThemeColor = 10
$ws.Range("I16").Interior.TintAndShade = 0

# J14 / J16 pick up the ordinary blue "in progress" marker used elsewhere
# in the sheet (same fill as J5, J7, ...).
$ws.Range("J5").Copy()
$ws.Range("J14").PasteSpecial(-4122)
$ws.Range("J16").PasteSpecial(-4122)

# H7 loses its old red highlight and becomes a plain grey cell like the
# other finished-week bars in that row (E8:H8 grey style).
$ws.Range("E8").Copy()
$ws.Range("H7").PasteSpecial(-4122)

# K7 / L7 and J17 / K17 (the old J16/K16, now shifted down) are cleared
# back to the "no fill, but explicitly formatted" style used by cells such
# as G5/H5.
$ws.Range("G5").Copy()
$ws.Range("K7").PasteSpecial(-4122)
$ws.Range("L7").PasteSpecial(-4122)
$ws.Range("J17").PasteSpecial(-4122)
$ws.Range("K17").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. View tweaks: zoom to 150% and move the active selection.
# ---------------------------------------------------------------------------
$wb.Windows.Item(1).Zoom = 150
$ws.Range("L13").Select()
